# "Added last minute updates" - update the topic placeholder paragraph:
#   - give the paragraph a 5-twip-space paragraph border on all 4 sides
#   - increase the left indent from 120 -> 225 twips
#   - rename the placeholder ID from ...topic_10... to ...501_3...
#   - drop the now-unneeded trailing " " run

$d = $word.ActiveDocument

# The placeholder/ID line is the first paragraph in the document.
$p1 = $d.Paragraphs.Item(1)

# --- remove the trailing whitespace-only run -------------------------
# It is the last run in the paragraph: a single space immediately before
# the paragraph mark, so it sits right before the paragraph's end.
$pEnd = $p1.Range.End
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# --- rename the ID placeholder text -----------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Find.Execute("**ID__AFFARS_5332_topic_10__ID**", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "**ID__AFFARS_5332_501_3__ID**", 2)

# --- paragraph formatting: indent + border ----------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
